# Scheduled runner update: refresh Marketboard profit calcs (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 794.9231
$ws.Range("I41").Value = 95
$ws.Range("J41").Value = 922.1818
$ws.Range("K41").Value = 95
$ws.Range("L41").Value = 922.1818
$ws.Range("M41").Value = 345
$ws.Range("N41").Value = -1802.1818

$ws.Range("H70").Value = 1280.1
$ws.Range("I70").Value = 1230.2
$ws.Range("J70").Value = 1330
$ws.Range("K70").Value = 3690.6
$ws.Range("L70").Value = 3990
$ws.Range("M70").Value = -3420.6
$ws.Range("N70").Value = -4530

$ws.Range("H73").Value = 1280.1
$ws.Range("I73").Value = 1230.2
$ws.Range("J73").Value = 1330
$ws.Range("K73").Value = 3690.6
$ws.Range("L73").Value = 3990
$ws.Range("M73").Value = -2754.6
$ws.Range("N73").Value = -5862

$ws.Range("H86").Value = 6405.316
$ws.Range("I86").Value = 1238.25
$ws.Range("J86").Value = 15263.143
$ws.Range("K86").Value = 1238.25
$ws.Range("L86").Value = 15263.143
$ws.Range("M86").Value = -115.25
$ws.Range("N86").Value = -17509.143

$ws.Range("H89").Value = 6405.316
$ws.Range("I89").Value = 1238.25
$ws.Range("J89").Value = 15263.143
$ws.Range("K89").Value = 6191.25
$ws.Range("L89").Value = 76315.715
$ws.Range("M89").Value = -575.25
$ws.Range("N89").Value = -87547.715

$ws.Range("H107").Value = 692
$ws.Range("I107").Value = 635.7273
$ws.Range("J107").Value = 1001.5
$ws.Range("K107").Value = 635.7273
$ws.Range("L107").Value = 1001.5
$ws.Range("M107").Value = 1284.2727
$ws.Range("N107").Value = -4841.5

$ws.Range("H129").Value = 176364.98
$ws.Range("J129").Value = 201024.4
$ws.Range("L129").Value = 603073.2
$ws.Range("N129").Value = -613073.2

$ws.Range("H132").Value = 2508.0476
$ws.Range("I132").Value = 2848
$ws.Range("K132").Value = 8544
$ws.Range("M132").Value = -6014

$ws.Range("H138").Value = 30305756
$ws.Range("I138").Value = 71430456
$ws.Range("J138").Value = 3343
$ws.Range("K138").Value = 214291368
$ws.Range("L138").Value = 10029
$ws.Range("M138").Value = -214286228
$ws.Range("N138").Value = -20309

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4828.5444
$ws.Range("I32").Value = 4060.8357
$ws.Range("K32").Value = 4060.8357
$ws.Range("M32").Value = -3773.8357

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H45").Value = 2463.6843
$ws.Range("I45").Value = 1776.7693
$ws.Range("K45").Value = 1776.7693
$ws.Range("M45").Value = -1399.7693

$ws.Range("H122").Value = 1578.5526
$ws.Range("I122").Value = 1446.1034
$ws.Range("J122").Value = 2005.3334
$ws.Range("K122").Value = 4338.3102
$ws.Range("L122").Value = 6016.0002
$ws.Range("M122").Value = -1888.3102
$ws.Range("N122").Value = -10916.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2730.4
$ws.Range("I20").Value = 2965.8333
$ws.Range("J20").Value = 2377.25
$ws.Range("K20").Value = 2965.8333
$ws.Range("L20").Value = 2377.25
$ws.Range("M20").Value = -2718.8333
$ws.Range("N20").Value = -2871.25

$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3318.3635
$ws.Range("I31").Value = 2887.6428
$ws.Range("J31").Value = 3519.3667
$ws.Range("K31").Value = 2887.6428
$ws.Range("L31").Value = 3519.3667
$ws.Range("M31").Value = -2592.6428
$ws.Range("N31").Value = -4109.3667

$ws.Range("H34").Value = 3318.3635
$ws.Range("I34").Value = 2887.6428
$ws.Range("J34").Value = 3519.3667
$ws.Range("K34").Value = 2887.6428
$ws.Range("L34").Value = 3519.3667
$ws.Range("M34").Value = -2685.6428
$ws.Range("N34").Value = -3923.3667

$ws.Range("H133").Value = 39150
$ws.Range("J133").Value = 39150
$ws.Range("L133").Value = 39150
$ws.Range("N133").Value = -44210

$ws.Range("H134").Value = 1335.8334
$ws.Range("I134").Value = 1233
$ws.Range("K134").Value = 3699
$ws.Range("M134").Value = -1164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 95.333336
$ws.Range("I40").Value = 95.333336
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 381.333344
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -312.333344
$ws.Range("N40").ClearContents()

$ws.Range("H56").Value = 8714.286
$ws.Range("I56").Value = 8714.286
$ws.Range("K56").Value = 8714.286
$ws.Range("M56").Value = -8184.286

$ws.Range("H117").Value = 1198.6666
$ws.Range("J117").Value = 1151.8
$ws.Range("L117").Value = 3455.4
$ws.Range("N117").Value = -10339.4

$ws.Range("H121").Value = 1193.375
$ws.Range("J121").Value = 1415
$ws.Range("L121").Value = 4245
$ws.Range("N121").Value = -6865

$ws.Range("H131").Value = 737.45
$ws.Range("J131").Value = 738.3333
$ws.Range("L131").Value = 2214.9999
$ws.Range("N131").Value = -12294.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 18367.4
$ws.Range("J113").Value = 2956.5
$ws.Range("L113").Value = 2956.5
$ws.Range("N113").Value = -7296.5

$ws.Range("H132").Value = 27474.25
$ws.Range("I132").Value = 2431.1428
$ws.Range("J132").Value = 85908.164
$ws.Range("K132").Value = 7293.428400000001
$ws.Range("L132").Value = 257724.492
$ws.Range("M132").Value = -4763.428400000001
$ws.Range("N132").Value = -262784.492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1282.8572
$ws.Range("I22").Value = 1746.5714
$ws.Range("J22").Value = 819.1429000000001
$ws.Range("K22").Value = 1746.5714
$ws.Range("L22").Value = 819.1429000000001
$ws.Range("M22").Value = -1451.5714
$ws.Range("N22").Value = -1409.1429

$ws.Range("H27").Value = 1282.8572
$ws.Range("I27").Value = 1746.5714
$ws.Range("J27").Value = 819.1429000000001
$ws.Range("K27").Value = 1746.5714
$ws.Range("L27").Value = 819.1429000000001
$ws.Range("M27").Value = -1639.5714
$ws.Range("N27").Value = -1033.1429

$ws.Range("H46").Value = 1119.6666
$ws.Range("I46").Value = 1010.19446
$ws.Range("K46").Value = 1010.19446
$ws.Range("M46").Value = -822.19446

$ws.Range("H93").Value = 1778.8
$ws.Range("I93").Value = 1830
$ws.Range("J93").Value = 1702
$ws.Range("K93").Value = 1830
$ws.Range("L93").Value = 1702
$ws.Range("M93").Value = -582
$ws.Range("N93").Value = -4198

$ws.Range("H116").Value = 34995
$ws.Range("J116").Value = 34995
$ws.Range("L116").Value = 34995
$ws.Range("N116").Value = -44173

$ws.Range("H118").Value = 500015000
$ws.Range("J118").Value = 500015000
$ws.Range("L118").Value = 500015000
$ws.Range("N118").Value = -500018314

$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -49676
